$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quadro.1")

# Remove the blank spacer row (old row 25) - this shifts the rows that
# contained "Total Valores Acrescentados...", "Impostos sobre os Produtos",
# "IVA", "Direitos de Importação", "Outros Impostos...", and
# "Produto Interno Bruto" each up by one row, and shrinks the used range
# from A1:AJ32 down to A1:AJ31.
$ws.Rows.Item(25).Delete()

# Reflect the new active selection on the bottom-right (frozen) pane, which
# now points at the relocated "Total Valores Acrescentados..." row.
$ws.Range("A25").Select()
